$d = $word.ActiveDocument

# 1. "(Forthcoming) " -> "2024" + tab
$null = $d.Content.Find.Execute("(Forthcoming) ", $true, $false, $false, $false, $false, $true, 1, $false, "2024^t", 2)

# 2. ", in the " -> ", "  (drop "in the" before the journal name)
$null = $d.Content.Find.Execute(", in the ", $true, $false, $false, $false, $false, $true, 1, $false, ", ", 2)

# 3. Replace the trailing "." that ends the paragraph (kept in its own, non-italic run)
#    with the expanded citation + line break + doi.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Carbon Offsets and Shifting Harms*") {
        $rng = $p.Range
        $endRng = $d.Range($rng.End - 2, $rng.End - 1)
        $null = $endRng.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, " 17(1), 234" + [char]0x2013 + "255. Special issue on intra- and interpersonal dilemmas in ethics and rational choice. With a critical comment by Kian Mintz-Woo and a response from me.^ldoi: 10.23941/ejpe.v17i1.790", 2)
        break
    }
}
